$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mercenaries")
$ws.Activate()

# Fixed mercenary icon root: icon files moved out of the icons/ "iconMerc*"
# naming into the plain "merc*" naming used elsewhere.
$ws.Range("N2").Value = "mercFootMan.png"
$ws.Range("N3").Value = "mercCleric.png"
$ws.Range("N4").Value = "mercCommander.png"
$ws.Range("N5").Value = "mercMage.png"
$ws.Range("N6").Value = "mercAssassin.png"
$ws.Range("N7").Value = "mercWarlock.png"

# Move the workbook window / active cell selection the way it was left
# after the edit session.
$ws.Range("N6").Select()

$win = $excel.ActiveWindow
$win.Left = 12090
